$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '82.106.89'
$ws.Range("E2").Value = '  +3.12%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.195.84'
$ws.Range("E3").Value = '  -0.30%  '

$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.00'
$ws.Range("E5").Value = '  +5.00%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '627.54'
$ws.Range("E6").Value = '  -1.15%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.288'
$ws.Range("E7").Value = '  +21.21%  '

$ws.Range("E8").Value = '  +0.13%  '

$ws.Range("E9").Value = '  +0.18%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.193.12'
$ws.Range("E10").Value = '  -0.50%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.591'
$ws.Range("E11").Value = '  -0.11%  '

$ws.Range("E12").Value = '  +13.11%  '

$ws.Range("E13").Value = '  -0.33%  '

$ws.Range("B14").Value = 'Toncoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.31'
$ws.Range("E14").Value = '  -3.64%  '

$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.777.70'
$ws.Range("E15").Value = '  -0.09%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '31.74'
$ws.Range("E16").Value = '  -0.78%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '81.895.00'
$ws.Range("E17").Value = '  +3.06%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.189.66'
$ws.Range("E18").Value = '  +0.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.20'
$ws.Range("E19").Value = '  +5.58%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.04'
$ws.Range("E20").Value = '  -3.59%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '435.99'
$ws.Range("E21").Value = '  +1.75%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.94'
$ws.Range("E22").Value = '  -2.99%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.11'
$ws.Range("E23").Value = '  +0.11%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.27'
$ws.Range("E24").Value = '  +5.72%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.27'
$ws.Range("E25").Value = '  +10.53%  '

$ws.Range("B26").Value = 'WrappedeETH'
$ws.Range("C26").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.348.25'
$ws.Range("E26").Value = '  -0.31%  '

$ws.Range("B27").Value = 'Litecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '76.77'
$ws.Range("E27").Value = '  -0.44%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.94'
$ws.Range("E28").Value = '  -3.16%  '

$ws.Range("E29").Value = '  -0.01%  '

$ws.Range("E30").Value = '  +3.44%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '593.28'
$ws.Range("E31").Value = '  +12.97%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.06'
$ws.Range("E32").Value = '  +0.15%  '

$ws.Range("E33").Value = '  +0.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.51'
$ws.Range("E34").Value = '  +1.12%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.155'
$ws.Range("E35").Value = '  +7.85%  '

$ws.Range("B36").Value = 'PancakeSwap'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.00'
$ws.Range("E36").Value = '  +0.05%  '

$ws.Range("B37").Value = 'Cronos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.138'
$ws.Range("E37").Value = '  +15.57%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '22.81'
$ws.Range("E38").Value = '  -0.50%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  +0.09%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.17'
$ws.Range("E40").Value = '  +11.36%  '

$ws.Range("E41").Value = '  +0.07%  '

$ws.Range("E42").Value = '  +14.08%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.08'
$ws.Range("E43").Value = '  +22.69%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '20.81'
$ws.Range("E44").Value = '  +3.84%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '160.76'
$ws.Range("E45").Value = '  -2.47%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '187.86'
$ws.Range("E47").Value = '  -2.32%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '44.41'
$ws.Range("E48").Value = '  +3.02%  '

$ws.Range("E49").Value = '  +0.35%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.777'
$ws.Range("E50").Value = '  -5.19%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '26.28'
$ws.Range("E51").Value = '  +1.71%  '
